$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion summary text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.79 = 6513.43 pesos`n✅ 6513.43 pesos = 1.78 = 945.19 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update the rate table values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 559
$ws2.Range("O10").Value = 3641.01
$ws2.Range("N12").Value = 3653
$ws2.Range("O12").Value = 530.1
